$d = $word.ActiveDocument

$replacements = @(
    @("69×87=6003", "48×98=4704"),
    @("78×19=1482", "85×67=5695"),
    @("82×82=6724", "29×71=2059"),
    @("14×14=196", "14×74=1036"),
    @("75×50=3750", "24×88=2112"),
    @("71×70=4970", "41×47=1927"),
    @("21×56=1176", "45×38=1710"),
    @("26×77=2002", "23×73=1679"),
    @("64×17=1088", "51×90=4590"),
    @("59×21=1239", "23×35=805"),
    @("57×86=4902", "39×47=1833"),
    @("32×43=1376", "74×47=3478"),
    @("19×79=1501", "86×11=946"),
    @("57×85=4845", "55×49=2695"),
    @("54×66=3564", "52×57=2964"),
    @("11×91=1001", "37×50=1850"),
    @("48×18=864", "93×43=3999"),
    @("23×22=506", "46×92=4232"),
    @("82×93=7626", "70×43=3010"),
    @("20×21=420", "49×23=1127"),
    @("84×73=6132", "40×59=2360"),
    @("45×47=2115", "46×32=1472"),
    @("76×69=5244", "47×35=1645"),
    @("33×69=2277", "72×37=2664"),
    @("13×11=143", "45×66=2970")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
